# Appends a new "Gathering data" section (heading + one bulleted item)
# after the "Making and breaking connections ..." bullet at the end of
# the "Larger model" section, preceded by a blank paragraph - matching
# the target diff.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Locate the paragraph that currently ends the document's content so the
# new material lands in the right place regardless of incidental shifts.
$searchRange = $d.Content
$found = $searchRange.Find.Execute(
    "Making and breaking connections is established for eight-friends but might need to have different dynamics in a larger model",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorPara = $searchRange.Paragraphs(1)

# 1) A blank paragraph with no style/numbering of its own.
$anchorPara.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Range.InsertXML("<w:p $wNs/>")

# 2) A Heading1 paragraph titled "Gathering data".
$p1 = $d.Paragraphs.Last
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Range.InsertXML("<w:p $wNs><w:pPr><w:pStyle w:val=`"Heading1`"/></w:pPr><w:r><w:t>Gathering data</w:t></w:r></w:p>")

# 3) A list-paragraph bullet (same numbering list used elsewhere, numId 2).
$p2 = $d.Paragraphs.Last
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Range.InsertXML("<w:p $wNs><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"2`"/></w:numPr></w:pPr><w:r><w:t>Use model_step to take the dimensionality and output it to a csv?</w:t></w:r></w:p>")
